$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(55,1).Value = "Auspost"
$ws.Cells.Item(55,2).Value = "Credit Card Account"
$ws.Cells.Item(56,1).Value = "Food Dairy"
$ws.Cells.Item(56,2).Value = "Business Account"
$ws.Cells.Item(57,1).Value = "Shopify"
$ws.Cells.Item(57,2).Value = "Credit Card Account"
$ws.Cells.Item(58,1).Value = "Adobe"
$ws.Cells.Item(58,2).Value = "Credit Card Account"
$ws.Cells.Item(59,1).Value = "Ikea"
$ws.Cells.Item(59,2).Value = "Credit Card Account"
$ws.Cells.Item(60,1).Value = "FYI"
$ws.Cells.Item(60,2).Value = "Credit Card Account"
$ws.Cells.Item(61,1).Value = "KP Lawyers"
$ws.Cells.Item(61,2).Value = "Business Account"
$ws.Cells.Item(62,1).Value = "Food Safety Inspection"
$ws.Cells.Item(62,2).Value = "Business Account"
$ws.Cells.Item(63,1).Value = "sinosmart"
$ws.Cells.Item(63,2).Value = "Business Account"
$ws.Cells.Item(64,1).Value = "legal"
$ws.Cells.Item(64,2).Value = "Business Account"
$ws.Cells.Item(65,1).Value = "HW Accounting"
$ws.Cells.Item(65,2).Value = "Business Account"
$ws.Cells.Item(66,1).Value = "daiwa"
$ws.Cells.Item(66,2).Value = "Business Account"
$ws.Cells.Item(67,1).Value = "de toni"
$ws.Cells.Item(67,2).Value = "Business Account"
$ws.Cells.Item(68,1).Value = "pomona"
$ws.Cells.Item(68,2).Value = "Business Account"
$ws.Cells.Item(69,1).Value = "Madhouse"
$ws.Cells.Item(69,2).Value = "Business Account"
$ws.Cells.Item(70,1).Value = "Tulip"
$ws.Cells.Item(70,2).Value = "Business Account"
$ws.Cells.Item(71,1).Value = "MF"
$ws.Cells.Item(71,2).Value = "Business Account"
$ws.Cells.Item(72,1).Value = "Fresh"
$ws.Cells.Item(72,2).Value = "Business Account"
$ws.Cells.Item(73,1).Value = "YCC"
$ws.Cells.Item(73,2).Value = "Business Account"
$ws.Cells.Item(74,1).Value = "Munja"
$ws.Cells.Item(74,2).Value = "Business Account"
$ws.Cells.Item(75,1).Value = "igeno"
$ws.Cells.Item(75,2).Value = "Home Loan Account"
$ws.Cells.Item(76,1).Value = "nova"
$ws.Cells.Item(76,2).Value = "Home Loan Account"
$ws.Cells.Item(77,1).Value = "bombora"
$ws.Cells.Item(77,2).Value = "Credit Card Account"
$ws.Cells.Item(78,1).Value = "1Password"
$ws.Cells.Item(78,2).Value = "Credit Card Account"
$ws.Cells.Item(79,1).Value = "Microsoft"
$ws.Cells.Item(79,2).Value = "Credit Card Account"
$ws.Cells.Item(80,1).Value = "Expedia Hotel"
$ws.Cells.Item(80,2).Value = "Credit Card Account"
$ws.Cells.Item(81,1).Value = "Agoda"
$ws.Cells.Item(81,2).Value = "Credit Card Account"
$ws.Cells.Item(82,1).Value = "Kmall"
$ws.Cells.Item(82,2).Value = "Credit Card Account"
$ws.Cells.Item(83,1).Value = "WL Filter"
$ws.Cells.Item(83,2).Value = "Credit Card Account"
$ws.Cells.Item(84,1).Value = "Cleanaway"
$ws.Cells.Item(84,2).Value = "Business Account"
$ws.Cells.Item(85,1).Value = "Victoria Basement"
$ws.Cells.Item(85,2).Value = "Credit Card Account"
$ws.Cells.Item(86,1).Value = "Premier Tazze"
$ws.Cells.Item(86,2).Value = "Credit Card Account"
$ws.Cells.Item(87,1).Value = "Coles"
$ws.Cells.Item(87,2).Value = "Credit Card Account"
$ws.Cells.Item(88,1).Value = "T2"
$ws.Cells.Item(88,2).Value = "Credit Card Account"
$ws.Cells.Item(89,1).Value = "New Yenyen"
$ws.Cells.Item(89,2).Value = "Credit Card Account"
$ws.Cells.Item(90,1).Value = "New renren"
$ws.Cells.Item(90,2).Value = "Credit Card Account"
$ws.Cells.Item(91,1).Value = "Crocs"
$ws.Cells.Item(91,2).Value = "Credit Card Account"
$ws.Cells.Item(92,1).Value = "Metro Petrol"
$ws.Cells.Item(92,2).Value = "Credit Card Account"
$ws.Cells.Item(93,1).Value = "Medco Petrol"
$ws.Cells.Item(93,2).Value = "Credit Card Account"